# samostatne fce pro typy prikladu
# Populate A1:A4 with the example's stored fields and style the header
# cell (A1) with a bold font, a thin box border and centered/top
# alignment - matching the "st12345 / john / doe / 5" record layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- values -----------------------------------------------------------
$ws.Range("A1").Value = "st12345"
$ws.Range("A2").Value = "john"
$ws.Range("A3").Value = "doe"

# --- A1 formatting: bold, thin box border, centered / top aligned -----
$a1 = $ws.Range("A1")
$a1.Font.Bold = $true
$a1.HorizontalAlignment = -4108   # xlCenter
$a1.VerticalAlignment = -4160     # xlTop
$a1.Borders.LineStyle = 1         # xlContinuous
$a1.Borders.Weight = 2            # xlThin

# --- A4: store "5" as text, not as a number ----------------------------
# Force text entry so the digit string is kept verbatim (not coerced to
# a numeric cell), then drop the now-unneeded Text number format so the
# cell itself stays unstyled, same as A2/A3.
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "5"
$ws.Range("A4").ClearFormats()
